$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new HGLG11 row is being inserted as row 2, which pushes the existing
# JSRE11 row down to row 3. Write the (unchanged) JSRE11 values into row 3
# first, then overwrite row 2 with the new HGLG11 values.
$ws.Range("A3").Value = "JSRE11"
$ws.Range("B3").Value = "74,03"
$ws.Range("C3").Value = "0,66"
$ws.Range("D3").Value = "7,78"
$ws.Range("E3").Value = "0,4600"
$ws.Range("F3").Value = "0,6806"
$ws.Range("G3").Value = "67,59"
$ws.Range("H3").Value = "15/05/2023"
$ws.Range("I3").Value = "0,4600"
$ws.Range("J3").Value = "0,6133"
$ws.Range("K3").Value = "75,00"
$ws.Range("L3").Value = "15/06/2023"

$ws.Range("A2").Value = "HGLG11"
$ws.Range("B2").Value = "158,05"
$ws.Range("C2").Value = "1,04"
$ws.Range("D2").Value = "10,44"
$ws.Range("E2").Value = "1,1000"
$ws.Range("F2").Value = "0,6770"
$ws.Range("G2").Value = "162,49"
$ws.Range("H2").Value = "15/05/2023"
$ws.Range("I2").Value = "1,1000"
$ws.Range("J2").Value = "0,6919"
$ws.Range("K2").Value = "158,98"
$ws.Range("L2").Value = "15/06/2023"

# Highlight the P/VP ("C") column: yellow for the new HGLG11 row, green for
# the JSRE11 row that got pushed down to row 3.
$ws.Range("C2").Interior.Color = 65535
$ws.Range("C3").Interior.Color = 65280
